$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.344867587089539
$ws.Range("B1").Value = 1.905002355575562
$ws.Range("C1").Value = 2.76426362991333
$ws.Range("D1").Value = 4.945603370666504
$ws.Range("E1").Value = 1.076334714889526
